$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the updated cryptos snapshot.
$changes = [ordered]@{
    "D2" = "42.830.13"
    "E2" = "  -5.17%  "
    "D3" = "2.207.45"
    "E3" = "  -6.53%  "
    "E4" = "  -0.03%  "
    "D5" = "315.22"
    "E5" = "  +1.07%  "
    "D6" = "98.74"
    "E6" = "  -8.35%  "
    "E7" = "  -7.04%  "
    "E8" = "  -0.13%  "
    "D9" = "0.559"
    "E9" = "  -8.33%  "
    "D10" = "36.75"
    "E10" = "  -10.04%  "
    "D11" = "53.95"
    "E11" = "  -3.11%  "
    "D12" = "0.0827"
    "E12" = "  -9.88%  "
    "D13" = "7.66"
    "E13" = "  -9.59%  "
    "E14" = "  -2.12%  "
    "B15" = "Polygon"
    "C15" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D15" = "0.856"
    "E15" = "  -12.19%  "
    "B16" = "WrappedliquidstakedEther2.0"
    "C16" = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
    "D16" = "2.542.25"
    "E16" = "  -6.48%  "
    "D17" = "14.21"
    "E17" = "  -6.67%  "
    "D18" = "2.206.41"
    "E18" = "  -6.83%  "
    "D19" = "42.720.86"
    "E19" = "  -5.34%  "
    "D20" = "14.39"
    "E20" = "  +0.73%  "
    "D21" = "0.0₃0959"
    "E21" = "  -9.68%  "
    "D22" = "6.41"
    "E22" = "  -10.87%  "
    "D23" = "65.12"
    "E23" = "  -11.14%  "
    "E24" = "  -10.26%  "
    "D25" = "235.36"
    "E25" = "  -9.22%  "
    "D26" = "2.12"
    "E26" = "  -8.22%  "
    "D27" = "0.998"
    "E27" = "  -0.41%  "
    "E28" = "  +1.58%  "
    "D29" = "9.97"
    "E29" = "  -9.72%  "
    "E30" = "  -4.21%  "
    "D31" = "6.28"
    "E31" = "  -12.74%  "
    "B32" = "EthereumClassic"
    "C32" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D32" = "20.51"
    "E32" = "  -8.14%  "
    "B33" = "Hedera"
    "C33" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D33" = "0.0882"
    "E33" = "  -9.72%  "
    "D34" = "34.17"
    "E34" = "  -8.28%  "
    "D35" = "154.80"
    "E35" = "  -7.70%  "
    "D36" = "2.79"
    "E36" = "  -6.80%  "
    "E37" = "  +7.35%  "
    "E38" = "  -7.05%  "
    "D39" = "1.91"
    "E39" = "  +6.76%  "
    "E40" = "  -6.88%  "
    "D41" = "4.40"
    "E41" = "  -5.76%  "
    "D42" = "3.77"
    "E42" = "  -4.46%  "
    "D43" = "0.0324"
    "E43" = "  -8.19%  "
    "D44" = "1.852.91"
    "E44" = "  +0.53%  "
    "E45" = "  +0.05%  "
    "D46" = "12.24"
    "E46" = "  -4.65%  "
    "D47" = "87.64"
    "E47" = "  -12.19%  "
    "E48" = "  -9.29%  "
    "D49" = "5.33"
    "E49" = "  -6.60%  "
    "D50" = "60.56"
    "E50" = "  -13.12%  "
    "D51" = "75.42"
    "E51" = "  -9.71%  "
}

foreach ($ref in $changes.Keys) {
    $value = $changes[$ref]
    $cell = $ws.Range($ref)
    # The source data stores every cell (prices, %-changes, names, links) as plain
    # text, even when it looks like a number (e.g. "315.22"). Force text formatting
    # first so Excel does not silently reinterpret it as a numeric value, then drop
    # back to the default "Normal" style so no stray formatting is left behind.
    $looksNumeric = $value -match "^[+-]?[0-9]*\.?[0-9]+$"
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
